# Add two new time-entries (rows) before the summary block, update one
# existing end-time, and let the sum/average rows shift down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the old "separator" row (row 71), which
# pushes the separator + 3 summary rows down from 71-74 to 73-76.
$ws.Rows("71:72").Insert()

# Correct the end time of the existing entry on row 70 (18:00 -> 16:30).
$ws.Range("E70").Value = 0.6875

# New entry: 2014-03-16, 19:10 - 20:20
$ws.Range("A71").Value = 2014
$ws.Range("B71").Value = 3
$ws.Range("C71").Value = 16
$ws.Range("D71").Value = 0.79861111111111116
$ws.Range("E71").Value = 0.84722222222222221
$ws.Range("F71").Formula = "=(E71-D71)*24*60"
$ws.Range("G71").Formula = "=F71/60"

# New entry: 2014-03-17, 09:25 - 11:00
$ws.Range("A72").Value = 2014
$ws.Range("B72").Value = 3
$ws.Range("C72").Value = 17
$ws.Range("D72").Value = 0.3923611111111111
$ws.Range("E72").Value = 0.45833333333333331
$ws.Range("F72").Formula = "=(E72-D72)*24*60"
$ws.Range("G72").Formula = "=F72/60"

# Update the current view selection to reflect where editing left off.
$ws.Range("J69").Select() | Out-Null
